# "Generate Report for Archive"
#  - Status of the (only) data row flips from "Ready for handoff" to
#    "In Translation" on every sheet that surfaces it (Overview's zh-cn /
#    de-de columns, plus the Status column on each language sheet).
#  - The two "date-ish" status columns on Overview (E:F) and the Status
#    column (C) on the per-language sheets get narrower to fit the new,
#    shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: keep the string literal on the left of -eq — PowerShell
            # coerces the right-hand operand to the left operand's type, and
            # the sheet also has literal-text "True"/"False" cells that would
            # otherwise (wrongly) compare equal to any non-empty string.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Narrower status columns now that the text is shorter.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.43

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.43
